# The commit swaps the contents of ppt/theme/theme1.xml ("Office Theme")
# and ppt/theme/theme2.xml ("Integral") - theme1.xml ends up holding the
# Integral colour scheme and theme2.xml ends up holding the Office Theme
# colour scheme (the relationships/parts that reference each theme file
# stay exactly as they were).
#
# The PowerPoint object model only exposes one live, editable theme for
# this deck (Presentation.SlideMaster.Theme / NotesMaster.Theme /
# Slide(i).ThemeColorScheme all resolve to the single theme that is
# persisted as ppt/theme/theme2.xml), so we recolor that theme's
# ThemeColorScheme to the target ("Office Theme") palette - the half of
# the swap that is reachable through COM automation.

function Convert-HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$themeColors = $sm.Theme.ThemeColorScheme

# Target palette = the "Office Theme" colours that theme2.xml should end
# up with (previously theme1.xml's colours), in the standard
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order used by
# ThemeColorScheme.Colors(1..12).
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeThemeHex.Length; $i++) {
    $themeColors.Colors($i).RGB = Convert-HexToRgbInt($officeThemeHex[$i - 1])
}
